$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 6 (year 2025) metrics
$ws.Range("C6").Value = 428
$ws.Range("E6").Value = 119
$ws.Range("G6").Value = 27.80373831775701
$ws.Range("H6").Value = 72.19626168224299
